# Actualización de planilla de avance.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# C2: was "en proceso" (shared string), now 80% progress value with percentage format.
$ws.Range("C2").Value = 0.8
$ws.Range("C2").NumberFormat = "0%"

# C3: new cell, 80% progress value with percentage format.
$ws.Range("C3").Value = 0.8
$ws.Range("C3").NumberFormat = "0%"

# C8: new cell, "en proceso" status text.
$ws.Range("C8").Value = "en proceso"

# Update the selected/active cell on the sheet from A3 to A21.
$ws.Range("A21").Select() | Out-Null
